$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 123.29412
$ws.Range("I33").Value = 99.875
$ws.Range("J33").Value = 498
$ws.Range("K33").Value = 99.875
$ws.Range("L33").Value = 498
$ws.Range("M33").Value = 129.125
$ws.Range("N33").Value = -956

$ws.Range("H55").Value = 137
$ws.Range("I55").Value = 137
$ws.Range("K55").Value = 137
$ws.Range("M55").Value = 77

$ws.Range("H92").Value = 361.57144
$ws.Range("I92").Value = 351.7
$ws.Range("J92").Value = 386.25
$ws.Range("K92").Value = 351.7
$ws.Range("L92").Value = 386.25
$ws.Range("M92").Value = 896.3
$ws.Range("N92").Value = -2882.25

$ws.Range("H96").Value = 862.55554
$ws.Range("I96").Value = 720.5
$ws.Range("J96").Value = 1999
$ws.Range("K96").Value = 2161.5
$ws.Range("L96").Value = 5997
$ws.Range("M96").Value = -788.5
$ws.Range("N96").Value = -8743

$ws.Range("H97").Value = 298
$ws.Range("J97").Value = 298
$ws.Range("L97").Value = 894
$ws.Range("N97").Value = -1886

$ws.Range("H99").Value = 326.75
$ws.Range("J99").Value = 299
$ws.Range("L99").Value = 897
$ws.Range("N99").Value = -3893

$ws.Range("H100").Value = 1853.0834
$ws.Range("I100").Value = 2395.4
$ws.Range("J100").Value = 1465.7142
$ws.Range("K100").Value = 2395.4
$ws.Range("L100").Value = 1465.7142
$ws.Range("M100").Value = -1854.4
$ws.Range("N100").Value = -2547.7142

$ws.Range("H112").Value = 1166.875
$ws.Range("J112").Value = 1476.8518
$ws.Range("L112").Value = 4430.555399999999
$ws.Range("N112").Value = -6646.555399999999

$ws.Range("H113").Value = 2994.1667
$ws.Range("I113").Value = 2991.25
$ws.Range("K113").Value = 2991.25
$ws.Range("M113").Value = 262.75

$ws.Range("H115").Value = 7634.5
$ws.Range("I115").Value = 7439.4287
$ws.Range("J115").Value = 9000
$ws.Range("K115").Value = 22318.2861
$ws.Range("L115").Value = 27000
$ws.Range("M115").Value = -20751.2861
$ws.Range("N115").Value = -30134

$ws.Range("H116").Value = 5847.5
$ws.Range("I116").Value = 5487
$ws.Range("J116").Value = 6929
$ws.Range("K116").Value = 5487
$ws.Range("L116").Value = 6929
$ws.Range("M116").Value = -2045
$ws.Range("N116").Value = -13813

$ws.Range("H118").Value = 1593.3334
$ws.Range("I118").Value = 990
$ws.Range("J118").Value = 2800
$ws.Range("K118").Value = 2970
$ws.Range("L118").Value = 8400
$ws.Range("M118").Value = -1313
$ws.Range("N118").Value = -11714

$ws.Range("H127").Value = 5338
$ws.Range("I127").Value = 797
$ws.Range("J127").Value = 7608.5
$ws.Range("K127").Value = 2391
$ws.Range("L127").Value = 22825.5
$ws.Range("M127").Value = 2569
$ws.Range("N127").Value = -32745.5

$ws.Range("H137").Value = 3690.4614
$ws.Range("J137").Value = 5434.875
$ws.Range("L137").Value = 16304.625
$ws.Range("N137").Value = -21404.625

$ws.Range("H138").Value = 2741.1455
$ws.Range("I138").Value = 1746.8572
$ws.Range("K138").Value = 5240.571599999999
$ws.Range("M138").Value = -100.5715999999993

$ws.Range("H141").Value = 3888.5386
$ws.Range("I141").Value = 3795.9167
$ws.Range("K141").Value = 11387.7501
$ws.Range("M141").Value = -6207.750100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1324.6666
$ws.Range("I45").Value = 1324.6666
$ws.Range("K45").Value = 1324.6666
$ws.Range("M45").Value = -947.6666

$ws.Range("H61").Value = 1563.909
$ws.Range("I61").Value = 1138.3125
$ws.Range("J61").Value = 2698.8333
$ws.Range("K61").Value = 1138.3125
$ws.Range("L61").Value = 2698.8333
$ws.Range("M61").Value = -926.3125
$ws.Range("N61").Value = -3122.8333

$ws.Range("H74").Value = 3052.75
$ws.Range("I74").Value = 3052.75
$ws.Range("K74").Value = 3052.75
$ws.Range("M74").Value = -2178.75

$ws.Range("H77").Value = 3052.75
$ws.Range("I77").Value = 3052.75
$ws.Range("K77").Value = 15263.75
$ws.Range("M77").Value = -10895.75

$ws.Range("H136").Value = 1563.909
$ws.Range("I136").Value = 1138.3125
$ws.Range("J136").Value = 2698.8333
$ws.Range("K136").Value = 3414.9375
$ws.Range("L136").Value = 8096.499899999999
$ws.Range("M136").Value = -864.9375
$ws.Range("N136").Value = -13196.4999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4475
$ws.Range("J105").Value = 7500
$ws.Range("L105").Value = 7500
$ws.Range("N105").Value = -10994

$ws.Range("H107").Value = 1743.375
$ws.Range("I107").Value = 1621.6666
$ws.Range("K107").Value = 1621.6666
$ws.Range("M107").Value = 298.3334

$ws.Range("H134").Value = 1882.55
$ws.Range("I134").Value = 1697
$ws.Range("K134").Value = 5091
$ws.Range("M134").Value = -2556

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1105.4615
$ws.Range("I107").Value = 564.6667
$ws.Range("K107").Value = 564.6667
$ws.Range("M107").Value = 1355.3333

$ws.Range("H134").Value = 2249.6667
$ws.Range("I134").Value = 2207.3333
$ws.Range("J134").Value = 2397.8333
$ws.Range("K134").Value = 6621.999899999999
$ws.Range("L134").Value = 7193.499899999999
$ws.Range("M134").Value = -4086.999899999999
$ws.Range("N134").Value = -12263.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 172.8
$ws.Range("I2").Value = 221.66667
$ws.Range("J2").Value = 99.5
$ws.Range("K2").Value = 1330.00002
$ws.Range("L2").Value = 597
$ws.Range("M2").Value = -1217.00002
$ws.Range("N2").Value = -823

$ws.Range("H33").Value = 1872
$ws.Range("I33").Value = 300
$ws.Range("J33").Value = 2658
$ws.Range("K33").Value = 1800
$ws.Range("L33").Value = 15948
$ws.Range("M33").Value = -1517
$ws.Range("N33").Value = -16514

$ws.Range("H56").Value = 18140
$ws.Range("I56").Value = 18140
$ws.Range("K56").Value = 18140
$ws.Range("M56").Value = -17610

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2282.389
$ws.Range("I132").Value = 1170.8572
$ws.Range("K132").Value = 3512.5716
$ws.Range("M132").Value = -982.5715999999998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4145.75
$ws.Range("I40").Value = 4145.75
$ws.Range("K40").Value = 4145.75
$ws.Range("M40").Value = -4009.75

$ws.Range("H55").Value = 133.66667
$ws.Range("J55").Value = 200
$ws.Range("L55").Value = 200
$ws.Range("N55").Value = -546

$ws.Range("H132").Value = 3774.375
$ws.Range("I132").Value = 2400
$ws.Range("J132").Value = 4599
$ws.Range("K132").Value = 7200
$ws.Range("L132").Value = 13797
$ws.Range("M132").Value = -4670
$ws.Range("N132").Value = -18857

$ws.Range("H136").Value = 5799.4
$ws.Range("I136").Value = 5799.4
$ws.Range("K136").Value = 17398.2
$ws.Range("M136").Value = -14848.2
